$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date values from F2 and F3 while keeping their existing style/format
$ws.Range("F2").ClearContents()
$ws.Range("F3").ClearContents()

# Update the active selection to F10
$ws.Range("F10").Select()
